$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")

# Row 76
$ws.Range("H76").Value = 4440
$ws.Range("I76").Value = 4042.8572
$ws.Range("J76").Value = 5366.6665
$ws.Range("K76").Value = 4042.8572
$ws.Range("L76").Value = 5366.6665
$ws.Range("M76").Value = -3727.8572
$ws.Range("N76").Value = -5996.6665

# Row 79
$ws.Range("H79").Value = 4440
$ws.Range("I79").Value = 4042.8572
$ws.Range("J79").Value = 5366.6665
$ws.Range("K79").Value = 4042.8572
$ws.Range("L79").Value = 5366.6665
$ws.Range("M79").Value = -2950.8572
$ws.Range("N79").Value = -7550.6665

# Row 98
$ws.Range("H98").Value = 951.25
$ws.Range("I98").Value = 951.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 951.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 546.75
$ws.Range("N98").ClearContents()

# Row 122
$ws.Range("H122").Value = 951.25
$ws.Range("I122").Value = 951.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2853.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -403.75
$ws.Range("N122").ClearContents()

# Row 125
$ws.Range("H125").Value = 857.1429000000001
$ws.Range("I125").Value = 800
$ws.Range("K125").Value = 7200
$ws.Range("M125").Value = -4740

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 11348.891
$ws.Range("I32").Value = 3105.173
$ws.Range("J32").Value = 31761.904
$ws.Range("K32").Value = 3105.173
$ws.Range("L32").Value = 31761.904
$ws.Range("M32").Value = -2818.173
$ws.Range("N32").Value = -32335.904

# Row 94
$ws.Range("H94").Value = 23333.334
$ws.Range("J94").Value = 23333.334
$ws.Range("L94").Value = 23333.334
$ws.Range("N94").Value = -25135.334

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 1797.4419
$ws.Range("I20").Value = 1729.7037
$ws.Range("J20").Value = 1911.75
$ws.Range("K20").Value = 1729.7037
$ws.Range("L20").Value = 1911.75
$ws.Range("M20").Value = -1482.7037
$ws.Range("N20").Value = -2405.75

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 11140203
$ws.Range("I31").Value = 6932111.5
$ws.Range("J31").Value = 28573728
$ws.Range("K31").Value = 6932111.5
$ws.Range("L31").Value = 28573728
$ws.Range("M31").Value = -6931816.5
$ws.Range("N31").Value = -28574318

# Row 34
$ws.Range("H34").Value = 11140203
$ws.Range("I34").Value = 6932111.5
$ws.Range("J34").Value = 28573728
$ws.Range("K34").Value = 6932111.5
$ws.Range("L34").Value = 28573728
$ws.Range("M34").Value = -6931909.5
$ws.Range("N34").Value = -28574132

# Row 122
$ws.Range("H122").Value = 4445316.5
$ws.Range("I122").Value = 6061327
$ws.Range("J122").Value = 1287.5
$ws.Range("K122").Value = 18183981
$ws.Range("L122").Value = 3862.5
$ws.Range("M122").Value = -18181531
$ws.Range("N122").Value = -8762.5

# Row 132
$ws.Range("H132").Value = 2089.111
$ws.Range("I132").Value = 1398.8572
$ws.Range("J132").Value = 2528.3635
$ws.Range("K132").Value = 4196.571599999999
$ws.Range("L132").Value = 7585.0905
$ws.Range("M132").Value = -1666.571599999999
$ws.Range("N132").Value = -12645.0905

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")

# Row 122
$ws.Range("H122").Value = 886.8125
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 927.7857
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 8350.0713
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -13250.0713

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 2327.9722
$ws.Range("I102").Value = 2410.9
$ws.Range("J102").Value = 1913.3334
$ws.Range("K102").Value = 2410.9
$ws.Range("L102").Value = 1913.3334
$ws.Range("M102").Value = -788.9000000000001
$ws.Range("N102").Value = -5157.3334

# Row 126
$ws.Range("H126").Value = 2134.2068
$ws.Range("I126").Value = 2249.375
$ws.Range("J126").Value = 1992.4615
$ws.Range("K126").Value = 6748.125
$ws.Range("L126").Value = 5977.3845
$ws.Range("M126").Value = -4278.125
$ws.Range("N126").Value = -10917.3845

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 2002
$ws.Range("I7").Value = 2002
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2002
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1890
$ws.Range("N7").ClearContents()

# Row 40
$ws.Range("H40").Value = 2651
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2651
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2651
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2923

# Row 122
$ws.Range("H122").Value = 5599.875
$ws.Range("I122").Value = 4799.8
$ws.Range("J122").Value = 6933.3335
$ws.Range("K122").Value = 14399.4
$ws.Range("L122").Value = 20800.0005
$ws.Range("M122").Value = -11949.4
$ws.Range("N122").Value = -25700.0005

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 125
$ws.Range("H125").Value = 47450
$ws.Range("J125").Value = 47450
$ws.Range("L125").Value = 47450
$ws.Range("N125").Value = -57290

# Row 126
$ws.Range("H126").Value = 2002
$ws.Range("I126").Value = 2002
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6006
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3536
$ws.Range("N126").ClearContents()

# Row 129
$ws.Range("H129").Value = 20000
$ws.Range("J129").Value = 20000
$ws.Range("L129").Value = 20000
$ws.Range("N129").Value = -30000

# Row 130
$ws.Range("H130").Value = 23357.25
$ws.Range("J130").Value = 23357.25
$ws.Range("L130").Value = 23357.25
$ws.Range("N130").Value = -33397.25

# Row 131
$ws.Range("H131").Value = 38775
$ws.Range("J131").Value = 38775
$ws.Range("L131").Value = 38775
$ws.Range("N131").Value = -48855

# Row 140
$ws.Range("H140").Value = 25000
$ws.Range("J140").Value = 25000
$ws.Range("L140").Value = 25000
$ws.Range("N140").Value = -35360

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 100001820
$ws.Range("I122").Value = 200001420
$ws.Range("J122").Value = 2230
$ws.Range("K122").Value = 600004260
$ws.Range("L122").Value = 6690
$ws.Range("M122").Value = -600001810
$ws.Range("N122").Value = -11590

# Row 123
$ws.Range("H123").Value = 46389.91
$ws.Range("J123").Value = 46389.91
$ws.Range("L123").Value = 46389.91
$ws.Range("N123").Value = -56189.91

# Row 124
$ws.Range("H124").Value = 60805.8
$ws.Range("J124").Value = 60805.8
$ws.Range("L124").Value = 60805.8
$ws.Range("N124").Value = -70625.8

# Row 125
$ws.Range("H125").Value = 120000
$ws.Range("J125").Value = 120000
$ws.Range("L125").Value = 120000
$ws.Range("N125").Value = -129840

# Row 127
$ws.Range("H127").Value = 17500
$ws.Range("J127").Value = 17500
$ws.Range("L127").Value = 17500
$ws.Range("N127").Value = -27420

# Row 128
$ws.Range("H128").Value = 28000
$ws.Range("J128").Value = 28000
$ws.Range("L128").Value = 28000
$ws.Range("N128").Value = -37960

Write-Output "edits applied"
